$d = $word.ActiveDocument

function Merge-Runs($searchText) {
    $rng = $d.Content
    $found = $rng.Find.Execute($searchText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if ($found) {
        # Re-assigning Range.Text collapses every run the range spans into a
        # single run that adopts the character formatting of the range's
        # first run. If the text we assign is byte-identical to what's
        # already there the engine treats it as a no-op and leaves the runs
        # split, so first swap in a throwaway placeholder to force a real
        # content change, then set the final (unchanged) text back.
        $rng.Text = "__TMP_MERGE_PLACEHOLDER__"
        $rng.Text = $searchText
    }
}

# Occurrence 1: "<id>p108r_1</id>" currently split across 3 runs
# (<id> / p108r_1 / </id>) -> merge into a single run.
Merge-Runs "<id>p108r_1</id>"

# Occurrence 2: "<id>p108r_2</id>" currently split across 4 runs
# (<id> / p108r_ / 2 / </id>) -> merge into a single run.
Merge-Runs "<id>p108r_2</id>"
